# Delete the "merchantId" column (column A) from the active sheet.
# This shifts merchantSiteId, gatewayName, userName, password, clientId, podId
# one column to the left (A:F), matching the frontend change where
# merchantId now defaults client-side instead of being listed in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).Delete()

$wb.Save()
